# Update progress ("ESTADO") values in the CRONOGRAMA table.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CRONOGRAMA")

# Row 3 - "Elaboracion del Documento de Gestion": 70% -> 100%
$ws.Range("G3").Value = 1

# Row 4 - "Elaboracion del Modelo del Negocio": 60% -> 100%
$ws.Range("G4").Value = 1

# Row 5 - "Elaboración de Lista de Requisitos": 0% -> 30%
$ws.Range("G5").Value = 0.3

# Leave the active selection on C5, matching the latest on-screen position.
$ws.Range("C5").Select()
